$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SCI")

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy cell formatting (number format/style) from column E into the new column D.
# Section-header rows (5, 6, 37, 79) only contain a label in column A/B and have
# no data columns, so they are excluded from the ranges below.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new/updated cell values (new column D plus a handful of cells whose
# figures were revised alongside the new column)
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 3190200
$ws.Range("D9").Value = 2429900
$ws.Range("E9").Value = 4745100
$ws.Range("F9").Value = 4709000
$ws.Range("D10").Value = 760300
$ws.Range("E10").Value = -1650100
$ws.Range("F10").Value = -1677900
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -5700
$ws.Range("E14").Value = 23100
$ws.Range("F14").Value = 66800
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 2569600
$ws.Range("E17").Value = 2524200
$ws.Range("F17").Value = 2540300
$ws.Range("D18").Value = 620500
$ws.Range("E18").Value = 570900
$ws.Range("F18").Value = 490800
$ws.Range("D20").Value = 2800
$ws.Range("E20").Value = -1500
$ws.Range("F20").Value = -2100
$ws.Range("D21").Value = 871800
$ws.Range("I21").Value = 711100
$ws.Range("D22").Value = 181600
$ws.Range("D23").Value = 441700
$ws.Range("D24").Value = 10300
$ws.Range("E24").Value = -400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 431500
$ws.Range("E26").Value = 400700
$ws.Range("D27").Value = 431100
$ws.Range("E27").Value = 400500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 16100
$ws.Range("E29").Value = 146200
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2800
$ws.Range("E32").Value = 1500
$ws.Range("F32").Value = 2100
$ws.Range("D33").Value = 447200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 447200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 177300
$ws.Range("D42").Value = 21500
$ws.Range("D43").Value = 82200
$ws.Range("D44").Value = 24900
$ws.Range("D45").Value = 25300
$ws.Range("D46").Value = 331200
$ws.Range("D47").Value = 5757800
$ws.Range("D48").Value = 3814800
$ws.Range("D49").Value = 2297700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 491700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 12693200
$ws.Range("D57").Value = 173400
$ws.Range("D58").Value = 86100
$ws.Range("D59").Value = 296100
$ws.Range("D60").Value = 555600
$ws.Range("D61").Value = 3532200
$ws.Range("D62").Value = 6963600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 11051300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 474300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1641900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 447200
$ws.Range("D83").Value = 248500
$ws.Range("I83").Value = 192400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 615800
$ws.Range("E89").Value = 503400
$ws.Range("F89").Value = 489000
$ws.Range("I89").Value = 385500
$ws.Range("D91").Value = -250100
$ws.Range("I91").Value = -112900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -414600
$ws.Range("E94").Value = -242900
$ws.Range("F94").Value = -221000
$ws.Range("D96").Value = -123800
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -329200
$ws.Range("F100").Value = -209500
$ws.Range("D101").Value = -5000
$ws.Range("I101").Value = -1000
$ws.Range("D102").Value = -133000
$ws.Range("E102").Value = 129100
$ws.Range("I102").Value = 52800
